# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Add a new "ODI Bowling Extra" sheet (mirrors the existing "ODI Batting Extra"
#    sheet) with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.
# 2) Tidy up "ODI Batting Extra" row 3 (match 3943) which had a few stray
#    blank cells (C3:E3) that should not be present.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clean up "ODI Batting Extra": row 3 should not carry the empty
#    C3/D3/E3 placeholder cells any more.
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("C3:E3").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet, right after "ODI Batting Extra".
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row, formatted the same way as the other "Extra" sheet headers.
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $bowlingExtra.Cells.Item(1, $col)
    $cell.Value2 = $headers[$col - 1]
}
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# Row data: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL.
# $null entries mean "leave the cell blank" (matches the source rows that
# have no recorded value for that column).
$rows = @(
    @("3939", "0", "40.00%"),
    @("3943", "0", $null),
    @("3944", "0", "30.00%"),
    @("3972", "0", "10.00%"),
    @("4114", "1", $null),
    @("4178", $null, $null),
    @("4200", "1", "30.00%"),
    @("4201", "0", $null),
    @("4204", "0", $null),
    @("4376", "0", "10.00%"),
    @("4460", "0", "30.00%"),
    @("4586", $null, $null),
    @("4590", $null, $null),
    @("4592", "1", "20.00%"),
    @("4634", $null, $null),
    @("4638", "0", "30.00%"),
    @("4641", "1", $null),
    @("4686", $null, $null),
    @("4688", "0", "40.00%"),
    @("4690", $null, $null)
)

$rowIndex = 2
foreach ($record in $rows) {
    $matchCode = $record[0]
    $maidenOvers = $record[1]
    $percentWickets = $record[2]

    $codeCell = $bowlingExtra.Cells.Item($rowIndex, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value2 = $matchCode

    if ($maidenOvers -ne $null) {
        $maidenCell = $bowlingExtra.Cells.Item($rowIndex, 2)
        $maidenCell.NumberFormat = "@"
        $maidenCell.Value2 = $maidenOvers
    }

    if ($percentWickets -ne $null) {
        $percentCell = $bowlingExtra.Cells.Item($rowIndex, 3)
        $percentCell.NumberFormat = "@"
        $percentCell.Value2 = $percentWickets
    }

    $rowIndex++
}
